# Add a new 4th-level bullet ("Tail Segments of the Snake") at the end of
# the body placeholder on the "Goals" slide (slide 2), right after the
# existing "Boarders of game board " bullet.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(2)
$textRange = $shape.TextFrame.TextRange

# Insert a new paragraph after all existing text; it inherits the
# indent level (lvl=3) of the preceding paragraph, matching PowerPoint's
# normal "press Enter at end of bullet, type text" behavior.
$newRange = $textRange.InsertAfter("`rTail Segments of the Snake")
